$wb = $excel.ActiveWorkbook

# Add a new worksheet after the existing sheets, named "TestCase02"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "TestCase02"

# Fill in header row
$newSheet.Range("A1").Value = "#"
$newSheet.Range("B1").Value = "Product"

# Fill in data rows
$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = "boss"

$newSheet.Range("A3").Value = 2
$newSheet.Range("B3").Value = "Surya"

$newSheet.Range("A4").Value = 3
$newSheet.Range("B4").Value = "Busch"

# Apply header style (copied from TestCase01's header row) and data style
$ws1 = $wb.Worksheets.Item("TestCase01")
$ws1.Range("A1").Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$newSheet.Range("A2:B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set column B width
$newSheet.Columns.Item(2).ColumnWidth = 17.5703125

# Update selection on sheet2 (TestCase01_1)
$ws2 = $wb.Worksheets.Item("TestCase01_1")
$ws2.Activate()
$ws2.Range("A1:D4").Select()

# Select B4 on new sheet, make it the active sheet (last, so it stays the active tab)
$newSheet.Activate()
$newSheet.Range("B4").Select()
